$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 119, shifting existing rows 119-165 down to 120-166.
$ws.Rows.Item(119).Insert()

# Populate the newly inserted row 119 with the new record's data.
$ws.Cells.Item(119, 1).Value = 8
$ws.Cells.Item(119, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(119, 3).Value = "Coquimbo"
$ws.Cells.Item(119, 4).Value = 44839
$ws.Cells.Item(119, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(119, 5).Value = 4
$ws.Cells.Item(119, 6).Value = 100112001
$ws.Cells.Item(119, 7).Value = "Berenjena"
$ws.Cells.Item(119, 8).Value = "Sin especificar"
$ws.Cells.Item(119, 9).Value = "Primera"
$ws.Cells.Item(119, 10).Value = 540
$ws.Cells.Item(119, 11).Value = 10800
$ws.Cells.Item(119, 12).Value = 11000
$ws.Cells.Item(119, 13).Value = 10900
$ws.Cells.Item(119, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(119, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(119, 16).Value = 272
$ws.Cells.Item(119, 17).Value = 40
$ws.Cells.Item(119, 18).Value = "Hortaliza"
